# Corrige o status (nomes x siglas) dos times no campeonato:
# "America-MG" e "Atletico-MG" estavam com as siglas trocadas.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Troca os nomes dos times nas linhas 2 e 4 (as siglas CAM/AME ja estavam certas)
$ws.Range("A2").Value = "Atlético-MG"
$ws.Range("A4").Value = "América-MG"

# Alinha a coluna A:B a esquerda (estava "general")
$ws.Range("A:B").HorizontalAlignment = -4131
